# Updates the cryptocurrency price/volume table on Sheet1 (columns
# B=Coin, C=Link, D=Price, E=Volume(1h)) to reflect the latest scrape.
# Most rows only get refreshed Price/Volume figures; rows 37-38 swap
# which coin (Hedera vs WEMIXToken) occupies which rank, and row 51
# replaces Algorand with RocketPoolETH entirely.
#
# Price values that look numeric (e.g. "308.42") are written with a
# text ("@") number format first so Excel stores them as text, matching
# the source data's inline-string cells instead of coercing to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "45.873.78"
$ws.Cells.Item(2, 5).Value = "  -0.62%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.591.70"
$ws.Cells.Item(3, 5).Value = "  -0.03%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.16%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "308.42"
$ws.Cells.Item(5, 5).Value = "  +0.19%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "98.60"
$ws.Cells.Item(6, 5).Value = "  -1.55%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.28%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.15%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.577"
$ws.Cells.Item(9, 5).Value = "  -0.40%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "38.51"
$ws.Cells.Item(10, 5).Value = "  +0.08%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "54.18"
$ws.Cells.Item(11, 5).Value = "  -1.01%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0839"
$ws.Cells.Item(12, 5).Value = "  +0.13%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  -2.51%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "2.992.67"
$ws.Cells.Item(14, 5).Value = "  +0.30%  "

# Row 15
$ws.Cells.Item(15, 5).Value = "  +1.64%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "2.596.61"
$ws.Cells.Item(16, 5).Value = "  +0.20%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.908"
$ws.Cells.Item(17, 5).Value = "  +0.75%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "14.72"
$ws.Cells.Item(18, 5).Value = "  -0.58%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "46.066.57"
$ws.Cells.Item(19, 5).Value = "  -0.41%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  +0.23%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.68"
$ws.Cells.Item(21, 5).Value = "  +0.61%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "12.54"
$ws.Cells.Item(22, 5).Value = "  -2.62%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "291.32"
$ws.Cells.Item(23, 5).Value = "  +14.83%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "73.03"
$ws.Cells.Item(24, 5).Value = "  +2.70%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "3.01"
$ws.Cells.Item(25, 5).Value = "  -0.36%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.24"
$ws.Cells.Item(26, 5).Value = "  +1.32%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "29.40"
$ws.Cells.Item(27, 5).Value = "  +4.61%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.999"
$ws.Cells.Item(28, 5).Value = "  -0.13%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "4.06"
$ws.Cells.Item(29, 5).Value = "  +0.98%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "10.68"
$ws.Cells.Item(30, 5).Value = "  +2.37%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "38.60"
$ws.Cells.Item(31, 5).Value = "  -2.31%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  -3.34%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "6.25"
$ws.Cells.Item(33, 5).Value = "  +3.48%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "3.62"
$ws.Cells.Item(34, 5).Value = "  -2.76%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "157.91"
$ws.Cells.Item(35, 5).Value = "  +3.62%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.21"
$ws.Cells.Item(36, 5).Value = "  -2.75%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "WEMIXToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.82"
$ws.Cells.Item(37, 5).Value = "  -2.63%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "Hedera"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.0836"
$ws.Cells.Item(38, 5).Value = "  +0.75%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +3.92%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  +0.58%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "15.59"
$ws.Cells.Item(41, 5).Value = "  -2.04%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.0327"
$ws.Cells.Item(42, 5).Value = "  +1.60%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "3.52"
$ws.Cells.Item(43, 5).Value = "  -2.31%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "3.98"
$ws.Cells.Item(44, 5).Value = "  -4.73%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "21.09"
$ws.Cells.Item(45, 5).Value = "  +3.39%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "2.109.87"
$ws.Cells.Item(46, 5).Value = "  +3.04%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.00"
$ws.Cells.Item(47, 5).Value = "  +0.07%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "95.25"
$ws.Cells.Item(48, 5).Value = "  +4.69%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "9.24"
$ws.Cells.Item(49, 5).Value = "  -0.25%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "108.52"
$ws.Cells.Item(50, 5).Value = "  -0.56%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "RocketPoolETH"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(51, 4).Value = "2.846.70"
$ws.Cells.Item(51, 5).Value = "  +0.03%  "
